# The tag group "3DD.003DE8629C" (16 detection rows, previously occupying
# worksheet rows 1115-1130) was removed from the export entirely. Deleting
# the entire-row range shifts every row below it up by 16, which is exactly
# what the target OOXML shows: the groups for tags "3DD.003E1189A4" and
# "3DD.003E16C1E6" move from rows 1131-1153 up to rows 1115-1137, and the
# sheet's used-range dimension shrinks from A1:P1153 to A1:P1137.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1115:1130").Delete()
